$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted price strings (not numbers) in the source data,
# so force text format before assigning to prevent Excel auto-converting
# numeric-looking strings (e.g. "237.92") into real numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.184.98"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.853.92"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.92"
$ws.Range("E5").Value = "  -1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6872"
$ws.Range("E6").Value = "  -4.54%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07664"
$ws.Range("E8").Value = "  +2.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3045"
$ws.Range("E9").Value = "  -3.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.19"
$ws.Range("E10").Value = "  -5.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08100"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.879.21"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7241"
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.196"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.37"
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.189.59"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.742"
$ws.Range("E17").Value = "  -4.29%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007799"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "233.69"
$ws.Range("E20").Value = "  -5.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.103.23"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.426"
$ws.Range("E24").Value = "  -3.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.76"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.961"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("E27").Value = "  -4.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.01"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.960"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.507"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.490"
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.011"
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05168"
$ws.Range("E34").Value = "  -4.01%  "
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7051"
$ws.Range("E36").Value = "  -3.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.024"
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.676"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01852"
$ws.Range("E39").Value = "  -3.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.683"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9361"
$ws.Range("E41").Value = "  +5.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.111.67"
$ws.Range("E42").Value = "  +6.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.956"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4281"
$ws.Range("E44").Value = "  -3.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.23"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.90"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.784"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.000.07"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.174"
$ws.Range("E50").Value = "  -4.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.994"
$ws.Range("E51").Value = "  -6.08%  "
